# informe del 16-11-2020 al 21-11-2020
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INFORME OCTUBRE Y NOVIEMBRE")
$ws.Activate()

# ------------------------------------------------------------------
# Row 29: the previous week's entry (row 31) moves up one slot in the
# log, so row 29 picks up that same content/height (165pt).
# ------------------------------------------------------------------
$ws.Range("A29").Value2 = "Modificacion el codigo en sql y C# de empresa, sucursal, vista usuario, vista empresa"
$ws.Range("B29").Value2 = "Modificar codigo y diseño en  sql y C# de formulario crear planilla y relaciones."
$ws.Range("C29").Value2 = "Modificacion codigo y diseño formulario documento, cargo, usuario, tambien se agrego nuevas columna en usurio y alterar todo el procedimiento en sql"
$ws.Range("D29").Value2 = "Tabla Mantenimiento de conceptos para la planilla en sql y su diseño en C# aun falta su codigo para el proceso en C#."
$ws.Range("E29").Value2 = "Campos en planilla(otros reintegros, prestacion alimentaria, recargo consumo, grati extra. Fiestas y navidad)"
$ws.Range("F29").Value2 = "Dias Vacaciones (separar los dias con validaciones en tabla  no subsidios no laborados, para su calculo respectivo. Calculo montos por trabajador(ONP, Y COMISIONES, SEGUROS, APORTES  DE AFP)"
$ws.Range("H29").Value2 = "CARLOS MEZA"
$ws.Rows.Item(29).RowHeight = 165

# ------------------------------------------------------------------
# Row 30: the date header row -- shift forward one week
# (16-nov-2020 .. 22-nov-2020)
# ------------------------------------------------------------------
$ws.Range("A30").Value2 = 44151
$ws.Range("B30").Value2 = 44152
$ws.Range("C30").Value2 = 44153
$ws.Range("D30").Value2 = 44154
$ws.Range("E30").Value2 = 44155
$ws.Range("F30").Value2 = 44156
$ws.Range("G30").Value2 = 44157

# ------------------------------------------------------------------
# Row 31: brand new content for the new week
# ------------------------------------------------------------------
$ws.Range("A31").Value2 = "nombre a los Controles CHECKBOX, para su validacion  con diseño completado"
$ws.Range("B31").Value2 = "Validar los controles  CHECKBOX para al hacer cuando hacemos click  ocultemos  los campos de la planilla, metodo para cada campo"
$ws.Range("C31").Value2 = "Procedimiento almacenado Registrar o Actualizar segun caso en SQL, Metodos set and get en C# para comunicar los parametros, clases(Dconcepto, Rconcepto, Nconcepto) "

$ws.Range("D31").Value2 = "Procedo registrar y pruebas tabla Conceptos en C#, error faltal al sincronizar  en el proyecto(se tuvo que hacer los cambios ya avazados desde el lunes"
$redRun = $ws.Range("D31").Characters(52, 100)
$redRun.Font.Color = 255
$redRun.Font.Name = "Calibri"
$redRun.Font.Size = 11

$ws.Range("E31").Value2 = "se modifico algoritmo calculo(ya que fue alterado al agreagar nuevos campos a la planilla), se agrego los nuevos campos a la tabla en la Base de datos. Se calculo aporte a O.N.P Y A.F.P , por cada trabajador"
$ws.Range("F31").Value2 = "formulario crear planilla y planilla calculo(se unieron y muestra los datos necesarios al seleccionar por periodo ok."
$ws.Range("H31").Value2 = "CARLOS MEZA"

# ------------------------------------------------------------------
# Row 32 was a leftover blank row under the table; remove it now that
# row 31 is the last week logged.
# ------------------------------------------------------------------
$ws.Rows.Item(32).Delete()

# ------------------------------------------------------------------
# Column H (DEVELOPER) no longer needs to be as wide
# ------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 23.25

# ------------------------------------------------------------------
# Leave the selection on the last edited cell
# ------------------------------------------------------------------
$ws.Range("F31").Select()

Write-Host "edit applied"
